$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: a record with numeric cells (PIM-9640 regression test fixture).
$ws.Range("A4").Value = "designer"
$ws.Range("B4").Value = 12345
$ws.Range("C4").Value = 12345
$ws.Range("D4").Value = 12345
$ws.Range("E4").Value = "Test with numeric values"
$ws.Range("G4").Value = 12345

# The whole new row (including the still-empty cells up to column K) picks up
# an explicit font-applied style, distinguishing it from the default style
# used by the header/data rows above.
$ws.Range("A4:K4").Font.Color = 0

# Selection ends up parked on G4, matching the edited workbook.
$null = $ws.Range("G4").Select()
